# Fill in Jesse's Week 5 Task Summary Sheet and Activity Log Summary Sheet
# with the local copies of his logs (per commit message: "copied over
# local copies of logs").

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# TASK SUMMARY SHEET
# ---------------------------------------------------------------------
$taskSheet = $wb.Worksheets.Item("TASK SUMMARY SHEET")

# Header: name + week number
$taskSheet.Range("C1").Value = "Jesse Hare"
$taskSheet.Range("E1").Value = 5

# Task rows (Stage / Task / Estimated Work Remaining at Start of Week /
# Hours Spent this Week / New Estimate)
$taskSheet.Range("A3").Value = "Project Build"
$taskSheet.Range("B3").Value = "Removing bugs with the searching algorithm and display"
$taskSheet.Range("C3").Value = 3
$taskSheet.Range("D3").Value = 8
$taskSheet.Range("E3").Value = 0

$taskSheet.Range("A4").Value = "Project Build"
$taskSheet.Range("B4").Value = "Converting db to an in-memory db"
$taskSheet.Range("C4").Value = 2
$taskSheet.Range("D4").Value = 1
$taskSheet.Range("E4").Value = 0

$taskSheet.Range("A5").Value = "Project Build"
$taskSheet.Range("B5").Value = "refactoring and optimisation"
$taskSheet.Range("C5").Value = 5
$taskSheet.Range("D5").Value = 3
$taskSheet.Range("E5").Value = 3

$taskSheet.Range("A6").Value = "Project Build"
$taskSheet.Range("B6").Value = "restructuring of code to make more modular for if new features needed"
$taskSheet.Range("C6").Value = 5
$taskSheet.Range("D6").Value = 6
$taskSheet.Range("E6").Value = 0

$taskSheet.Range("A7").Value = "Project Build"
$taskSheet.Range("B7").Value = "testing with dummy csv file, invalid input files and user input"
$taskSheet.Range("C7").Value = 5
$taskSheet.Range("D7").Value = 2
$taskSheet.Range("E7").Value = 0

# ---------------------------------------------------------------------
# ACTIVITY LOG SUMMARY SHEET
# ---------------------------------------------------------------------
$summarySheet = $wb.Worksheets.Item("ACTIVITY LOG SUMMARY SHEET")

$summarySheet.Range("D1").Value = "Jesse Hare"

$summarySheet.Range("A4").Value = "Project Build"
$summarySheet.Range("B4").Value = 18
$summarySheet.Range("C4").Value = 2
